$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Marking" row (row 11): points per right answer and per wrong answer
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Update the "Total" row (row 12): total marks and the score/max text
$ws.Range("B12").Value = 36
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "34 / 112"
